# Add a "plannedTarget" efficiency column (column J), pulled from the
# planning database, alongside the existing daily-plan columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, matching the formatting used by the rest of row 1.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "plannedTarget"

# Planned-target efficiency values fetched from the database for the
# first two dates (shift A and shift B rows each).
$ws.Range("J2").Value = 135
$ws.Range("J3").Value = 135
$ws.Range("J4").Value = 135
$ws.Range("J5").Value = 135

# Leave the selection where the user was last working.
$ws.Range("J8").Select()
